$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the missing Friday hours for Jake (week 2) ---
$ws.Range("H13").Value = 2

# --- New column J: labels + values for Project 3 person-hour tracking ---

# J6: "Total Person Hours for Project 3:" (underlined label)
$ws.Range("J6").Value = "Total Person Hours for Project 3:"
$ws.Range("J6").Font.Underline = $true

# J7: array-entered total of both weekly tables, highlighted yellow, left-aligned
$ws.Range("J7").FormulaArray = "=SUM(B3:H7+B9:H13)"
$ws.Range("J7").Interior.Color = 65535
$ws.Range("J7").HorizontalAlignment = -4131

# J2: "Estimate of Person Hours for Project 3:" (underlined label)
$ws.Range("J2").Value = "Estimate of Person Hours for Project 3:"
$ws.Range("J2").Font.Underline = $true

# J3: the estimate value, highlighted yellow, left-aligned (reuse J7's
# format via copy/paste-special so the same style index is shared instead
# of minting a redundant one)
$ws.Range("J3").Value = 120
$ws.Range("J7").Copy()
$ws.Range("J3").PasteSpecial(-4122)

# Widen column J to fit the label text
$ws.Columns("J:J").ColumnWidth = 29.6

# Move the active selection, as in the saved file
$ws.Range("J18").Select()
